# Append the latest daily Covid-19 data point (25/5/2020) to the
# "Tabela1" table on the active worksheet, mirroring the row-76 entry
# produced by the nightly data-refresh bot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adding a ListRow expands the table's range (and the worksheet
# dimension) the same way Excel does when a table is grown by one row.
$tbl = $ws.ListObjects.Item("Tabela1")
$newRow = $tbl.ListRows.Add()

# --- Column A: Date -------------------------------------------------
# Set formatting *before* the value so the number isn't coerced to
# text by the (default, text-oriented) format of a brand-new cell.
$a = $ws.Range("A76")
$a.NumberFormat = "d/\ m/\ yyyy;@"
$a.Font.Name = "Calibri Light"
$a.Font.Size = 10
$a.HorizontalAlignment = -4152
$a.VerticalAlignment = -4160
$a.Locked = $true
$a.FormulaHidden = $false
$a.Value = 43976

# --- Column B: Tested (all) -----------------------------------------
$b = $ws.Range("B76")
$b.NumberFormat = "#,##0"
$b.Font.Name = "Calibri Light"
$b.Font.Size = 10
$b.HorizontalAlignment = -4152
$b.VerticalAlignment = -4107
$b.Value = 75770

# --- Columns C:J: remaining daily counters ---------------------------
$rest = $ws.Range("C76:J76")
$rest.Font.Name = "Calibri Light"
$rest.Font.Size = 10
$rest.HorizontalAlignment = -4152
$rest.VerticalAlignment = -4107

$ws.Range("C76").Value = 754
$ws.Range("D76").Value = 1469
$ws.Range("E76").Value = 0
$ws.Range("F76").Value = 9
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 6
$ws.Range("I76").Value = 108
$ws.Range("J76").Value = 1

# Match the saved selection left behind by the bot's run.
$ws.Range("A76:J76").Select()
